$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header fixes -----------------------------------------------------
# "Nota Cieências" had a typo -> "Nota Ciências"
$ws.Range("E1").Value = "Nota Ciências"
# "Avaliação de Matematica" -> "Avaliação Matematica"
$ws.Range("F1").Value = "Avaliação Matematica"

# --- Extend the table with two new evaluation columns -----------------
# Pre-format the new columns (G:H) the same way as the rest of the table
# (copy the plain/general-alignment format used throughout the sheet)
# before writing any values into them, so the new cells share the same
# cell style as the existing ones.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("G1:H6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("G1").Value = "Avaliação Português"
$ws.Range("H1").Value = "Avaliação Ciências"

# --- New "Aprovado/Reprovado" formula for Claudio's math evaluation ---
$ws.Range("F2").Formula = "=if C2 < 7'Reprovado' else 'Aprovado'"

# Column F now holds this much longer text, so it is resized accordingly.
$ws.Columns.Item(6).ColumnWidth = 68.43
